# Update the "想去人数" (want-to-go count) column F values on the
# "展览" and "全部类型" worksheets to reflect the refreshed scrape.

$wb = $excel.ActiveWorkbook

# Row -> new F value for the "展览" sheet
$sheet1Updates = @{
    2  = 6632
    4  = 416
    5  = 65
    6  = 9
    7  = 541
    8  = 97
    11 = 2
    13 = 393
    14 = 1283
    15 = 11
    16 = 3309
    18 = 211
    19 = 1939
    20 = 62
    22 = 127
}

# Row -> new F value for the "全部类型" sheet
$sheet4Updates = @{
    2  = 6632
    4  = 416
    5  = 65
    6  = 9
    8  = 541
    9  = 97
    12 = 2
    14 = 393
    15 = 1283
    16 = 11
    17 = 3309
    19 = 211
    20 = 1939
    21 = 62
    23 = 127
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
